{"js": "// Adds three new paragraphs after the existing \"Exemplo\" paragraph:\n//   1. A plain paragraph (no style) -> \"Se add mais um par\u00e1grafo...\"\n//   2. A paragraph using the custom \"EstiloExemploPython\" style -> \"Mas se passar o estilo...\"\n//   3. A paragraph using the built-in \"Heading5\" style -> \"Tamb\u00e9m \u00e9 poss\u00edvel utilizar...\"\n//\n// New paragraphs inserted via `insertParagraph`/`InsertLocation.after` inherit the\n// style of the paragraph they follow (same as pressing Enter in Word), which is not\n// what the target markup shows (the first new paragraph has no <w:pPr> at all, and\n// the following ones carry their own explicit style). To reproduce the exact\n// structure we build each paragraph as a small flat-OPC OOXML fragment and insert it\n// with `insertOoxml`, which inserts the markup verbatim instead of inheriting\n// formatting from the insertion point.\n\nfunction flatOpcDocument(bodyXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + bodyXml + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n  );\n}\n\nfunction paragraphXml(text, styleId) {\n  const pPr = styleId ? `<w:pPr><w:pStyle w:val=\"${styleId}\"/></w:pPr>` : \"\";\n  return `<w:p>${pPr}<w:r><w:t>${text}</w:t></w:r></w:p>`;\n}\n\nasync function appendParagraphOoxml(text, styleId) {\n  const paragraphs = context.document.body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n\n  const lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n  const insertionPoint = lastParagraph.getRange(Word.RangeLocation.end);\n  const ooxml = flatOpcDocument(paragraphXml(text, styleId));\n  insertionPoint.insertOoxml(ooxml, Word.InsertLocation.after);\n  await context.sync();\n}\n\nawait appendParagraphOoxml(\n  \"Se add mais um par\u00e1grafo ele vir\u00e1 sem a formata\u00e7\u00e3o.\"\n);\nawait appendParagraphOoxml(\n  \"Mas se passar o estilo, ele ir\u00e1 formatar\",\n  \"EstiloExemploPython\"\n);\nawait appendParagraphOoxml(\n  \"Tamb\u00e9m \u00e9 poss\u00edvel utilizar estilos padr\u00f5es do Word\",\n  \"Heading5\"\n);\n", "ps1": "# Adds three new paragraphs after the existing \"Exemplo\" paragraph:\n#   1. A plain paragraph (no style) -> \"Se add mais um par\u00e1grafo...\"\n#   2. A paragraph using the custom \"EstiloExemploPython\" style -> \"Mas se passar o estilo...\"\n#   3. A paragraph using the built-in \"Heading5\" style -> \"Tamb\u00e9m \u00e9 poss\u00edvel utilizar...\"\n#\n# Plain Range.InsertParagraphAfter() (i.e. pressing Enter at the end of the\n# document) makes the new paragraph inherit the style of the paragraph before it,\n# which does not match the target markup: the first new paragraph has no <w:pPr>\n# at all, while the following two carry their own explicit style. To reproduce the\n# exact structure, each new paragraph is built as a small flat-OPC WordOpenXML\n# fragment and injected with Range.InsertXML, which inserts the markup verbatim\n# instead of inheriting formatting from the insertion point.\n\n$d = $word.ActiveDocument\n\nfunction New-FlatOpcDocument([string]$BodyXml) {\n    return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" ' +\n        'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' + $BodyXml + '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>'\n}\n\nfunction New-ParagraphXml([string]$Text, [string]$StyleId) {\n    $pPr = \"\"\n    if ($StyleId) {\n        $pPr = '<w:pPr><w:pStyle w:val=\"' + $StyleId + '\"/></w:pPr>'\n    }\n    return '<w:p>' + $pPr + '<w:r><w:t>' + $Text + '</w:t></w:r></w:p>'\n}\n\nfunction Add-ParagraphAtEnd([string]$Text, [string]$StyleId) {\n    $endRange = $d.Range($d.Content.End, $d.Content.End)\n    $ooxml = New-FlatOpcDocument (New-ParagraphXml $Text $StyleId)\n    $endRange.InsertXML($ooxml)\n}\n\nAdd-ParagraphAtEnd \"Se add mais um par\u00e1grafo ele vir\u00e1 sem a formata\u00e7\u00e3o.\" $null\nAdd-ParagraphAtEnd \"Mas se passar o estilo, ele ir\u00e1 formatar\" \"EstiloExemploPython\"\nAdd-ParagraphAtEnd \"Tamb\u00e9m \u00e9 poss\u00edvel utilizar estilos padr\u00f5es do Word\" \"Heading5\"\n"}
